$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = "Sistemas"
$ws.Range("E6").Select()
